$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("04-10-2021", 0.12, 0.9,  0.41, 1.02, 1.6,  2.69, 1.35, 2.62),
    @("05-10-2021", 0.13, 1.54, 0.32, 0.93, 2.14, 1.94, 1.12, 2.26),
    @("06-10-2021", $null, 0.62, 0.38, 0.91, 2.67, 2.17, 1.03, 2.23)
)

$startRow = 191
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $val = $row[$c]
        if ($null -ne $val) {
            $cell = $ws.Cells.Item($r, $c + 1)
            if ($c -eq 0) {
                $cell.Value = "'" + $val
                $cell.Style = "Normal"
            } else {
                $cell.Value = $val
            }
        }
    }
}
